$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.538.99"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.377.63"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.377.14"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "3.948.78"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").Value = "3.367.07"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "61.586.76"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").Value = "3.511.29"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +7.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  -6.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0771"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.772"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "2.352.42"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("E51").Value = "  +0.70%  "
